$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Week 32 (row 21): Mon/Tue/Wed are now marked done, matching the green
# "completed" formatting already used for earlier weeks (e.g. row 9), and
# 3 working days are logged for that week.
$ws.Range("D9").Copy()
$ws.Range("D21:F21").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("I21").Value = 3
